# Update the "想去人数" (interest count) figures in column F for the
# rows that were refreshed in this data pull (两个工作表： 展览 / 全部类型）。
$wb = $excel.ActiveWorkbook

# Map: worksheet name -> @{ row = newValue }
$updates = @{
    "展览"   = @{ 4 = 1167; 8 = 278; 12 = 525; 15 = 13098; 17 = 8; 19 = 5379; 20 = 5552 }
    "全部类型" = @{ 4 = 1167; 24 = 278; 34 = 525; 37 = 13098; 39 = 8; 42 = 5379; 43 = 5552 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
